{"js": "// Add a new paragraph \"Guten Morgen\" at the end of the document body,\n// right after the existing \"------------------------------------\" line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Insert a new paragraph after the last existing paragraph (the dashed\n// separator line) and set its text to \"Guten Morgen\".\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Guten Morgen\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new paragraph \"Guten Morgen\" at the end of the document body,\n# right after the existing \"------------------------------------\" line.\n$d = $word.ActiveDocument\n\n# Grab the last paragraph currently in the document (the dashed separator\n# line) and insert a brand-new paragraph mark right after it.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n# The newly created paragraph is now the last one in the document; give it\n# the requested text.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Guten Morgen\"\n"}
